{"js": "// Recolor specific \"Steps to Perform\" / hypothesis bullet points in the\n// marketing campaign problem statement:\n//   - Three EDA-related bullets become red (FF0000).\n//   - Three hypothesis bullets become green (00B050).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nconst redTexts = new Set([\n  \"Create box plots and histograms to understand the distributions and outliers. Perform outlier treatment.\",\n  \"Use ordinal encoding and one hot encoding according to different types of categorical variables.\",\n  \"Create a heatmap to showcase the correlation between different pairs of variables.\"\n]);\n\nconst greenTexts = new Set([\n  \"Older people are not as tech-savvy and probably prefer shopping in-store.\",\n  \"Customers with kids probably have less time to visit a store and would prefer to shop online.\",\n  \"Other distribution channels may cannibalize sales at the store.\"\n]);\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (redTexts.has(text)) {\n    para.font.color = \"#FF0000\";\n  } else if (greenTexts.has(text)) {\n    para.font.color = \"#00B050\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Recolor specific \"Steps to Perform\" / hypothesis bullet points in the\n# marketing campaign problem statement:\n#   - Three EDA-related bullets become red (wdColorRed / RGB(255,0,0)).\n#   - Three hypothesis bullets become green (RGB(0,176,80) = 00B050).\n\n$d = $word.ActiveDocument\n\n# Word's Font.Color is a BGR-packed long (0x00BBGGRR), not RRGGBB.\n$wdColorRed = 0x0000FF      # RRGGBB FF0000\n$wdColorGreen = 0x50B000    # RRGGBB 00B050\n\n$redTexts = @(\n    \"Create box plots and histograms to understand the distributions and outliers. Perform outlier treatment.\",\n    \"Use ordinal encoding and one hot encoding according to different types of categorical variables.\",\n    \"Create a heatmap to showcase the correlation between different pairs of variables.\"\n)\n\n$greenTexts = @(\n    \"Older people are not as tech-savvy and probably prefer shopping in-store.\",\n    \"Customers with kids probably have less time to visit a store and would prefer to shop online.\",\n    \"Other distribution channels may cannibalize sales at the store.\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($redTexts -contains $text) {\n        $p.Range.Font.Color = $wdColorRed\n    } elseif ($greenTexts -contains $text) {\n        $p.Range.Font.Color = $wdColorGreen\n    }\n}\n"}
